$d = $word.ActiveDocument

# The names in this table use non-breaking spaces (U+00A0) between words,
# so build the search/replacement strings with that character to preserve
# the original formatting convention.
$nbsp = [char]0x00A0

$oldKichoE = "Kicho" + $nbsp + "E" + $nbsp + "Jang"
$newKichoI = "Kicho" + $nbsp + "I" + $nbsp + "Jang"

$oldKichoSamJan = "Kicho" + $nbsp + "Sam" + $nbsp + "Jan"
$newKichoSamJang = "Kicho" + $nbsp + "Sam" + $nbsp + "Jang"

$oldPalgueE = "Palgue" + $nbsp + "E" + $nbsp + "Jang"
$newPalgueI = "Palgue" + $nbsp + "I" + $nbsp + "Jang"

# Fix "Kicho E Jang" -> "Kicho I Jang"
$d.Content.Find.Execute($oldKichoE, $true, $true, $false, $false, $false,
                         $true, 1, $false, $newKichoI, 2)

# Fix typo "Kicho Sam Jan" -> "Kicho Sam Jang"
$d.Content.Find.Execute($oldKichoSamJan, $true, $true, $false, $false, $false,
                         $true, 1, $false, $newKichoSamJang, 2)

# Fix "Palgue E Jang" -> "Palgue I Jang"
$d.Content.Find.Execute($oldPalgueE, $true, $true, $false, $false, $false,
                         $true, 1, $false, $newPalgueI, 2)
